$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_12_3_0"
$ws.Range("B2").Value = 0.7163343551151435
$ws.Range("C2").Value = 0.2116952027380578
$ws.Range("D2").Value = 0.8214247572930771
$ws.Range("E2").Value = 0.7027398888239587
$ws.Range("F2").Value = 0.3139342963695526
$ws.Range("G2").Value = 0.2002371996641159
$ws.Range("H2").Value = 0.1182246208190918
$ws.Range("I2").Value = 0.1616432517766953

$ws.Range("A3").Value = "model_12_3_1"
$ws.Range("B3").Value = 0.7560738676548382
$ws.Range("C3").Value = -0.1896065430967537
$ws.Range("D3").Value = 0.5794930260491453
$ws.Range("E3").Value = 0.464886856554132
$ws.Range("F3").Value = 0.2699543535709381
$ws.Range("G3").Value = 0.3021718263626099
$ws.Range("H3").Value = 0.2783939838409424
$ws.Range("I3").Value = 0.2909823060035706

$ws.Range("A4").Value = "model_12_3_22"
$ws.Range("B4").Value = 0.7676871150971099
$ws.Range("C4").Value = 0.09538780734470997
$ws.Range("D4").Value = 0.5019516181356489
$ws.Range("E4").Value = 0.4909392965204002
$ws.Range("F4").Value = 0.2571019232273102
$ws.Range("G4").Value = 0.2297804057598114
$ws.Range("H4").Value = 0.3297297954559326
$ws.Range("I4").Value = 0.2768155932426453

$ws.Range("A5").Value = "model_12_3_21"
$ws.Range("B5").Value = 0.7676871150971099
$ws.Range("C5").Value = 0.09538780734470997
$ws.Range("D5").Value = 0.5019516181356489
$ws.Range("E5").Value = 0.4909392965204002
$ws.Range("F5").Value = 0.2571019232273102
$ws.Range("G5").Value = 0.2297804057598114
$ws.Range("H5").Value = 0.3297297954559326
$ws.Range("I5").Value = 0.2768155932426453

$ws.Range("A6").Value = "model_12_3_20"
$ws.Range("B6").Value = 0.7676871150971099
$ws.Range("C6").Value = 0.09538780734470997
$ws.Range("D6").Value = 0.5019516181356489
$ws.Range("E6").Value = 0.4909392965204002
$ws.Range("F6").Value = 0.2571019232273102
$ws.Range("G6").Value = 0.2297804057598114
$ws.Range("H6").Value = 0.3297297954559326
$ws.Range("I6").Value = 0.2768155932426453

$ws.Range("A7").Value = "model_12_3_19"
$ws.Range("B7").Value = 0.7676871150971099
$ws.Range("C7").Value = 0.09538780734470997
$ws.Range("D7").Value = 0.5019516181356489
$ws.Range("E7").Value = 0.4909392965204002
$ws.Range("F7").Value = 0.2571019232273102
$ws.Range("G7").Value = 0.2297804057598114
$ws.Range("H7").Value = 0.3297297954559326
$ws.Range("I7").Value = 0.2768155932426453

$ws.Range("A8").Value = "model_12_3_18"
$ws.Range("B8").Value = 0.7676871150971099
$ws.Range("C8").Value = 0.09538780734470997
$ws.Range("D8").Value = 0.5019516181356489
$ws.Range("E8").Value = 0.4909392965204002
$ws.Range("F8").Value = 0.2571019232273102
$ws.Range("G8").Value = 0.2297804057598114
$ws.Range("H8").Value = 0.3297297954559326
$ws.Range("I8").Value = 0.2768155932426453

$ws.Range("A9").Value = "model_12_3_17"
$ws.Range("B9").Value = 0.7676871150971099
$ws.Range("C9").Value = 0.09538780734470997
$ws.Range("D9").Value = 0.5019516181356489
$ws.Range("E9").Value = 0.4909392965204002
$ws.Range("F9").Value = 0.2571019232273102
$ws.Range("G9").Value = 0.2297804057598114
$ws.Range("H9").Value = 0.3297297954559326
$ws.Range("I9").Value = 0.2768155932426453

$ws.Range("A10").Value = "model_12_3_16"
$ws.Range("B10").Value = 0.7676871150971099
$ws.Range("C10").Value = 0.09538780734470997
$ws.Range("D10").Value = 0.5019516181356489
$ws.Range("E10").Value = 0.4909392965204002
$ws.Range("F10").Value = 0.2571019232273102
$ws.Range("G10").Value = 0.2297804057598114
$ws.Range("H10").Value = 0.3297297954559326
$ws.Range("I10").Value = 0.2768155932426453

$ws.Range("A11").Value = "model_12_3_15"
$ws.Range("B11").Value = 0.7676871150971099
$ws.Range("C11").Value = 0.09538780734470997
$ws.Range("D11").Value = 0.5019516181356489
$ws.Range("E11").Value = 0.4909392965204002
$ws.Range("F11").Value = 0.2571019232273102
$ws.Range("G11").Value = 0.2297804057598114
$ws.Range("H11").Value = 0.3297297954559326
$ws.Range("I11").Value = 0.2768155932426453

$ws.Range("A12").Value = "model_12_3_14"
$ws.Range("B12").Value = 0.7676871150971099
$ws.Range("C12").Value = 0.09538780734470997
$ws.Range("D12").Value = 0.5019516181356489
$ws.Range("E12").Value = 0.4909392965204002
$ws.Range("F12").Value = 0.2571019232273102
$ws.Range("G12").Value = 0.2297804057598114
$ws.Range("H12").Value = 0.3297297954559326
$ws.Range("I12").Value = 0.2768155932426453

$ws.Range("A13").Value = "model_12_3_13"
$ws.Range("B13").Value = 0.7676871150971099
$ws.Range("C13").Value = 0.09538780734470997
$ws.Range("D13").Value = 0.5019516181356489
$ws.Range("E13").Value = 0.4909392965204002
$ws.Range("F13").Value = 0.2571019232273102
$ws.Range("G13").Value = 0.2297804057598114
$ws.Range("H13").Value = 0.3297297954559326
$ws.Range("I13").Value = 0.2768155932426453

$ws.Range("A14").Value = "model_12_3_12"
$ws.Range("B14").Value = 0.7676871150971099
$ws.Range("C14").Value = 0.09538780734470997
$ws.Range("D14").Value = 0.5019516181356489
$ws.Range("E14").Value = 0.4909392965204002
$ws.Range("F14").Value = 0.2571019232273102
$ws.Range("G14").Value = 0.2297804057598114
$ws.Range("H14").Value = 0.3297297954559326
$ws.Range("I14").Value = 0.2768155932426453

$ws.Range("A15").Value = "model_12_3_11"
$ws.Range("B15").Value = 0.7676871150971099
$ws.Range("C15").Value = 0.09538780734470997
$ws.Range("D15").Value = 0.5019516181356489
$ws.Range("E15").Value = 0.4909392965204002
$ws.Range("F15").Value = 0.2571019232273102
$ws.Range("G15").Value = 0.2297804057598114
$ws.Range("H15").Value = 0.3297297954559326
$ws.Range("I15").Value = 0.2768155932426453

$ws.Range("A16").Value = "model_12_3_10"
$ws.Range("B16").Value = 0.7676871150971099
$ws.Range("C16").Value = 0.09538780734470997
$ws.Range("D16").Value = 0.5019516181356489
$ws.Range("E16").Value = 0.4909392965204002
$ws.Range("F16").Value = 0.2571019232273102
$ws.Range("G16").Value = 0.2297804057598114
$ws.Range("H16").Value = 0.3297297954559326
$ws.Range("I16").Value = 0.2768155932426453

$ws.Range("A17").Value = "model_12_3_9"
$ws.Range("B17").Value = 0.7676871150971099
$ws.Range("C17").Value = 0.09538780734470997
$ws.Range("D17").Value = 0.5019516181356489
$ws.Range("E17").Value = 0.4909392965204002
$ws.Range("F17").Value = 0.2571019232273102
$ws.Range("G17").Value = 0.2297804057598114
$ws.Range("H17").Value = 0.3297297954559326
$ws.Range("I17").Value = 0.2768155932426453

$ws.Range("A18").Value = "model_12_3_8"
$ws.Range("B18").Value = 0.7676871150971099
$ws.Range("C18").Value = 0.09538780734470997
$ws.Range("D18").Value = 0.5019516181356489
$ws.Range("E18").Value = 0.4909392965204002
$ws.Range("F18").Value = 0.2571019232273102
$ws.Range("G18").Value = 0.2297804057598114
$ws.Range("H18").Value = 0.3297297954559326
$ws.Range("I18").Value = 0.2768155932426453

$ws.Range("A19").Value = "model_12_3_7"
$ws.Range("B19").Value = 0.7676871150971099
$ws.Range("C19").Value = 0.09538780734470997
$ws.Range("D19").Value = 0.5019516181356489
$ws.Range("E19").Value = 0.4909392965204002
$ws.Range("F19").Value = 0.2571019232273102
$ws.Range("G19").Value = 0.2297804057598114
$ws.Range("H19").Value = 0.3297297954559326
$ws.Range("I19").Value = 0.2768155932426453

$ws.Range("A20").Value = "model_12_3_6"
$ws.Range("B20").Value = 0.7676871150971099
$ws.Range("C20").Value = 0.09538780734470997
$ws.Range("D20").Value = 0.5019516181356489
$ws.Range("E20").Value = 0.4909392965204002
$ws.Range("F20").Value = 0.2571019232273102
$ws.Range("G20").Value = 0.2297804057598114
$ws.Range("H20").Value = 0.3297297954559326
$ws.Range("I20").Value = 0.2768155932426453

$ws.Range("A21").Value = "model_12_3_5"
$ws.Range("B21").Value = 0.7676871150971099
$ws.Range("C21").Value = 0.09538780734470997
$ws.Range("D21").Value = 0.5019516181356489
$ws.Range("E21").Value = 0.4909392965204002
$ws.Range("F21").Value = 0.2571019232273102
$ws.Range("G21").Value = 0.2297804057598114
$ws.Range("H21").Value = 0.3297297954559326
$ws.Range("I21").Value = 0.2768155932426453

$ws.Range("A22").Value = "model_12_3_4"
$ws.Range("B22").Value = 0.7676871150971099
$ws.Range("C22").Value = 0.09538780734470997
$ws.Range("D22").Value = 0.5019516181356489
$ws.Range("E22").Value = 0.4909392965204002
$ws.Range("F22").Value = 0.2571019232273102
$ws.Range("G22").Value = 0.2297804057598114
$ws.Range("H22").Value = 0.3297297954559326
$ws.Range("I22").Value = 0.2768155932426453

$ws.Range("A23").Value = "model_12_3_3"
$ws.Range("B23").Value = 0.7676871150971099
$ws.Range("C23").Value = 0.09538780734470997
$ws.Range("D23").Value = 0.5019516181356489
$ws.Range("E23").Value = 0.4909392965204002
$ws.Range("F23").Value = 0.2571019232273102
$ws.Range("G23").Value = 0.2297804057598114
$ws.Range("H23").Value = 0.3297297954559326
$ws.Range("I23").Value = 0.2768155932426453

$ws.Range("A24").Value = "model_12_3_2"
$ws.Range("B24").Value = 0.7676871150971099
$ws.Range("C24").Value = 0.09538780734470997
$ws.Range("D24").Value = 0.5019516181356489
$ws.Range("E24").Value = 0.4909392965204002
$ws.Range("F24").Value = 0.2571019232273102
$ws.Range("G24").Value = 0.2297804057598114
$ws.Range("H24").Value = 0.3297297954559326
$ws.Range("I24").Value = 0.2768155932426453

$ws.Range("A25").Value = "model_12_3_23"
$ws.Range("B25").Value = 0.7676871150971099
$ws.Range("C25").Value = 0.09538780734470997
$ws.Range("D25").Value = 0.5019516181356489
$ws.Range("E25").Value = 0.4909392965204002
$ws.Range("F25").Value = 0.2571019232273102
$ws.Range("G25").Value = 0.2297804057598114
$ws.Range("H25").Value = 0.3297297954559326
$ws.Range("I25").Value = 0.2768155932426453

$ws.Range("A26").Value = "model_12_3_24"
$ws.Range("B26").Value = 0.7676871150971099
$ws.Range("C26").Value = 0.09538780734470997
$ws.Range("D26").Value = 0.5019516181356489
$ws.Range("E26").Value = 0.4909392965204002
$ws.Range("F26").Value = 0.2571019232273102
$ws.Range("G26").Value = 0.2297804057598114
$ws.Range("H26").Value = 0.3297297954559326
$ws.Range("I26").Value = 0.2768155932426453
